$wb = $excel.ActiveWorkbook

# Sheet: f1_score
$ws = $wb.Worksheets.Item("f1_score")
$ws.Range("B2").Value = "0.327 (0.321 Â± 0.005)"
$ws.Range("C2").Value = "0.489 (0.476 Â± 0.010)"
$ws.Range("D2").Value = "0.741 (0.730 Â± 0.012)"
$ws.Range("E2").Value = "0.741 (0.726 Â± 0.016)"
$ws.Range("F2").Value = "0.258 (0.238 Â± 0.019)"
$ws.Range("B3").Value = "0.379 (0.304 Â± 0.039)"
$ws.Range("C3").Value = "0.516 (0.475 Â± 0.018)"
$ws.Range("D3").Value = "0.738 (0.705 Â± 0.014)"
$ws.Range("E3").Value = "0.799 (0.742 Â± 0.026)"
$ws.Range("F3").Value = "0.229 (0.191 Â± 0.021)"
$ws.Range("B4").Value = "0.351 (0.258 Â± 0.037)"
$ws.Range("C4").Value = "0.487 (0.430 Â± 0.034)"
$ws.Range("D4").Value = "0.695 (0.663 Â± 0.023)"
$ws.Range("E4").Value = "0.716 (0.679 Â± 0.024)"
$ws.Range("F4").Value = "0.209 (0.173 Â± 0.018)"
$ws.Range("B5").Value = "0.349 (0.275 Â± 0.043)"
$ws.Range("C5").Value = "0.474 (0.227 Â± 0.148)"
$ws.Range("D5").Value = "0.746 (0.533 Â± 0.283)"
$ws.Range("E5").Value = "0.800 (0.625 Â± 0.151)"
$ws.Range("F5").Value = "0.217 (0.171 Â± 0.023)"
$ws.Range("B6").Value = "0.363 (0.303 Â± 0.032)"
$ws.Range("C6").Value = "0.533 (0.492 Â± 0.026)"
$ws.Range("D6").Value = "0.764 (0.718 Â± 0.015)"
$ws.Range("E6").Value = "0.796 (0.749 Â± 0.021)"
$ws.Range("F6").Value = "0.230 (0.198 Â± 0.018)"
$ws.Range("C7").Value = "0.469 (0.468 Â± 0.001)"
$ws.Range("B9").Value = "0.369 (0.316 Â± 0.040)"
$ws.Range("C9").Value = "0.496 (0.464 Â± 0.022)"
$ws.Range("D9").Value = "0.745 (0.717 Â± 0.016)"
$ws.Range("E9").Value = "0.772 (0.723 Â± 0.024)"
$ws.Range("F9").Value = "0.224 (0.157 Â± 0.050)"
$ws.Range("C10").Value = "0.456 (0.450 Â± 0.005)"
$ws.Range("B11").Value = "0.281 (0.166 Â± 0.070)"
$ws.Range("C11").Value = "0.370 (0.307 Â± 0.045)"
$ws.Range("D11").Value = "0.565 (0.489 Â± 0.047)"
$ws.Range("E11").Value = "0.607 (0.472 Â± 0.143)"
$ws.Range("F11").Value = "0.180 (0.106 Â± 0.044)"
$ws.Range("B12").Value = "0.294 (0.274 Â± 0.016)"
$ws.Range("C12").Value = "0.506 (0.463 Â± 0.028)"
$ws.Range("D12").Value = "0.663 (0.337 Â± 0.325)"
$ws.Range("E12").Value = "0.773 (0.773 Â± 0.000)"
$ws.Range("B13").Value = "0.094 (0.045 Â± 0.023)"
$ws.Range("C13").Value = "0.146 (0.099 Â± 0.033)"
$ws.Range("D13").Value = "0.308 (0.173 Â± 0.081)"
$ws.Range("E13").Value = "0.299 (0.235 Â± 0.040)"
$ws.Range("F13").Value = "0.018 (0.005 Â± 0.005)"
$ws.Range("B16").Value = "0.389 (0.313 Â± 0.038)"
$ws.Range("C16").Value = "0.495 (0.455 Â± 0.020)"
$ws.Range("D16").Value = "0.752 (0.713 Â± 0.018)"
$ws.Range("E16").Value = "0.781 (0.726 Â± 0.022)"
$ws.Range("F16").Value = "0.224 (0.191 Â± 0.017)"
$ws.Range("B17").Value = "0.385 (0.305 Â± 0.033)"
$ws.Range("C17").Value = "0.504 (0.460 Â± 0.022)"
$ws.Range("D17").Value = "0.734 (0.702 Â± 0.019)"
$ws.Range("E17").Value = "0.753 (0.704 Â± 0.024)"
$ws.Range("F17").Value = "0.213 (0.184 Â± 0.017)"

# Sheet: training_time
$ws = $wb.Worksheets.Item("training_time")
$ws.Range("B2").Value = "00:04:47 (00:05:32 Â± 00:00:48)"
$ws.Range("C2").Value = "00:04:59 (00:08:18 Â± 00:02:34)"
$ws.Range("D2").Value = "00:04:51 (00:08:30 Â± 00:04:46)"
$ws.Range("E2").Value = "00:04:54 (00:10:22 Â± 00:06:53)"
$ws.Range("F2").Value = "00:06:10 (00:27:23 Â± 00:13:09)"
$ws.Range("B3").Value = "00:00:33 (00:00:50 Â± 00:00:21)"
$ws.Range("C3").Value = "00:01:05 (00:01:22 Â± 00:00:11)"
$ws.Range("D3").Value = "00:00:44 (00:02:21 Â± 00:01:50)"
$ws.Range("E3").Value = "00:01:28 (00:03:38 Â± 00:01:36)"
$ws.Range("F3").Value = "00:05:03 (00:06:33 Â± 00:01:07)"
$ws.Range("B4").Value = "00:00:26 (00:00:35 Â± 00:00:08)"
$ws.Range("C4").Value = "00:00:43 (00:00:59 Â± 00:00:16)"
$ws.Range("D4").Value = "00:03:25 (00:03:56 Â± 00:00:38)"
$ws.Range("E4").Value = "00:01:16 (00:01:40 Â± 00:00:21)"
$ws.Range("F4").Value = "00:01:00 (00:01:21 Â± 00:00:19)"
$ws.Range("B5").Value = "00:05:05 (00:05:13 Â± 00:00:04)"
$ws.Range("C5").Value = "00:05:05 (00:05:13 Â± 00:00:03)"
$ws.Range("D5").Value = "00:05:08 (00:05:14 Â± 00:00:03)"
$ws.Range("E5").Value = "00:05:07 (00:05:13 Â± 00:00:04)"
$ws.Range("F5").Value = "00:05:08 (00:05:16 Â± 00:00:04)"
$ws.Range("B6").Value = "00:04:58 (00:05:03 Â± 00:00:02)"
$ws.Range("C6").Value = "00:04:56 (00:05:01 Â± 00:00:02)"
$ws.Range("D6").Value = "00:05:00 (00:05:06 Â± 00:00:04)"
$ws.Range("E6").Value = "00:04:56 (00:05:00 Â± 00:00:02)"
$ws.Range("F6").Value = "00:04:57 (00:05:00 Â± 00:00:03)"
$ws.Range("C7").Value = "00:05:03 (00:05:06 Â± 00:00:02)"
$ws.Range("B9").Value = "00:04:59 (00:05:01 Â± 00:00:02)"
$ws.Range("C9").Value = "00:05:00 (00:05:04 Â± 00:00:03)"
$ws.Range("D9").Value = "00:05:00 (00:05:03 Â± 00:00:04)"
$ws.Range("E9").Value = "00:05:01 (00:05:04 Â± 00:00:02)"
$ws.Range("F9").Value = "00:05:00 (00:05:07 Â± 00:00:15)"
$ws.Range("C10").Value = "00:04:29 (00:04:29 Â± 00:00:00)"
$ws.Range("B11").Value = "00:05:05 (00:05:06 Â± 00:00:00)"
$ws.Range("C11").Value = "00:05:06 (00:05:07 Â± 00:00:00)"
$ws.Range("D11").Value = "00:05:03 (00:05:07 Â± 00:00:00)"
$ws.Range("E11").Value = "00:05:05 (00:05:06 Â± 00:00:00)"
$ws.Range("F11").Value = "00:05:16 (00:05:54 Â± 00:00:26)"
$ws.Range("B12").Value = "00:01:12 (00:02:13 Â± 00:00:36)"
$ws.Range("C12").Value = "00:02:03 (00:05:11 Â± 00:02:05)"
$ws.Range("D12").Value = "00:03:25 (00:03:43 Â± 00:00:17)"
$ws.Range("E12").Value = "00:02:04 (00:02:04 Â± 00:00:00)"
$ws.Range("B13").Value = "00:00:02 (00:00:03 Â± 00:00:01)"
$ws.Range("C13").Value = "00:00:08 (00:00:09 Â± 00:00:00)"
$ws.Range("D13").Value = "00:00:11 (00:00:14 Â± 00:00:02)"
$ws.Range("E13").Value = "00:00:19 (00:00:19 Â± 00:00:00)"
$ws.Range("F13").Value = "00:00:07 (00:00:08 Â± 00:00:00)"
$ws.Range("B16").Value = "00:12:43 (01:21:38 Â± 00:46:22)"
$ws.Range("C16").Value = "00:28:15 (01:50:14 Â± 01:06:35)"
$ws.Range("D16").Value = "00:09:39 (03:14:13 Â± 01:29:18)"
$ws.Range("E16").Value = "00:19:58 (00:58:47 Â± 00:46:34)"
$ws.Range("F16").Value = "00:09:28 (00:09:48 Â± 00:00:11)"
$ws.Range("B17").Value = "00:05:01 (00:05:26 Â± 00:00:17)"
$ws.Range("C17").Value = "00:05:03 (00:05:51 Â± 00:00:26)"
$ws.Range("D17").Value = "00:05:05 (00:06:10 Â± 00:01:52)"
$ws.Range("E17").Value = "00:05:03 (00:05:56 Â± 00:00:35)"
$ws.Range("F17").Value = "00:05:18 (00:06:48 Â± 00:01:27)"

# Sheet: test_time
$ws = $wb.Worksheets.Item("test_time")
$ws.Range("B2").Value = "00:00:01 (00:00:05 Â± 00:00:02)"
$ws.Range("C2").Value = "00:00:01 (00:00:04 Â± 00:00:02)"
$ws.Range("D2").Value = "00:00:03 (00:00:10 Â± 00:00:05)"
$ws.Range("E2").Value = "00:00:03 (00:00:09 Â± 00:00:04)"
$ws.Range("F2").Value = "00:00:00 (00:00:02 Â± 00:00:01)"
$ws.Range("B3").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("C3").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("D3").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("E3").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("F3").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("B4").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("C4").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("D4").Value = "00:00:09 (00:00:09 Â± 00:00:00)"
$ws.Range("E4").Value = "00:00:01 (00:00:01 Â± 00:00:00)"
$ws.Range("F4").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("B5").Value = "00:00:01 (00:00:02 Â± 00:00:01)"
$ws.Range("C5").Value = "00:00:01 (00:00:02 Â± 00:00:02)"
$ws.Range("D5").Value = "00:00:01 (00:00:02 Â± 00:00:02)"
$ws.Range("E5").Value = "00:00:01 (00:00:02 Â± 00:00:02)"
$ws.Range("F5").Value = "00:00:01 (00:00:02 Â± 00:00:01)"
$ws.Range("B6").Value = "00:00:00 (00:00:02 Â± 00:00:00)"
$ws.Range("C6").Value = "00:00:00 (00:00:02 Â± 00:00:01)"
$ws.Range("D6").Value = "00:00:00 (00:00:04 Â± 00:00:02)"
$ws.Range("E6").Value = "00:00:01 (00:00:05 Â± 00:00:01)"
$ws.Range("F6").Value = "00:00:01 (00:00:03 Â± 00:00:02)"
$ws.Range("C7").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("B9").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("C9").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("D9").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("E9").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("F9").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("C10").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("B11").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("C11").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("D11").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("E11").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("F11").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("B12").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("C12").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("D12").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("E12").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("B13").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("C13").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("D13").Value = "00:00:00 (00:00:01 Â± 00:00:00)"
$ws.Range("E13").Value = "00:00:01 (00:00:01 Â± 00:00:00)"
$ws.Range("F13").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("B16").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("C16").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("D16").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("E16").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("F16").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("B17").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("C17").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("D17").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("E17").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("F17").Value = "00:00:00 (00:00:00 Â± 00:00:00)"

# Sheet: missing_runs
$ws = $wb.Worksheets.Item("missing_runs")
$ws.Range("B2").Value = "[]"
$ws.Range("C2").Value = "[]"
$ws.Range("D2").Value = "[]"
$ws.Range("E2").Value = "[]"
$ws.Range("F2").Value = "[]"
$ws.Range("B3").Value = "[]"
$ws.Range("C3").Value = "[]"
$ws.Range("D3").Value = "[]"
$ws.Range("E3").Value = "[]"
$ws.Range("F3").Value = "[]"
$ws.Range("B4").Value = "[]"
$ws.Range("C4").Value = "[]"
$ws.Range("D4").Value = "[]"
$ws.Range("E4").Value = "[]"
$ws.Range("F4").Value = "[]"
$ws.Range("B5").Value = "[]"
$ws.Range("C5").Value = "[]"
$ws.Range("D5").Value = "[]"
$ws.Range("E5").Value = "[]"
$ws.Range("F5").Value = "[]"
$ws.Range("B6").Value = "[]"
$ws.Range("C6").Value = "[]"
$ws.Range("D6").Value = "[]"
$ws.Range("E6").Value = "[]"
$ws.Range("F6").Value = "[]"
$ws.Range("C7").Value = "[2, 3, 5, 7, 11, 13, 17, 19, 23, 29, 41, 43, 47, 53, 59, 61, 67, 71]"
$ws.Range("B9").Value = "[]"
$ws.Range("C9").Value = "[]"
$ws.Range("D9").Value = "[]"
$ws.Range("E9").Value = "[]"
$ws.Range("F9").Value = "[]"
$ws.Range("C10").Value = "[2, 3, 5, 7, 11, 13, 17, 19, 23, 29, 41, 43, 47, 53, 59, 61, 67, 71]"
$ws.Range("B11").Value = "[]"
$ws.Range("C11").Value = "[]"
$ws.Range("D11").Value = "[]"
$ws.Range("E11").Value = "[]"
$ws.Range("F11").Value = "[]"
$ws.Range("B12").Value = "[2, 3, 7, 11, 17, 31, 37, 41, 43, 53, 67, 71]"
$ws.Range("C12").Value = "[2, 5, 7, 11, 19, 23, 31, 37, 41, 67, 71]"
$ws.Range("D12").Value = "[2, 5, 7, 11, 13, 17, 19, 29, 31, 37, 41, 43, 47, 53, 59, 61, 67, 71]"
$ws.Range("E12").Value = "[2, 3, 5, 7, 13, 17, 19, 23, 29, 31, 37, 41, 43, 47, 53, 59, 61, 67, 71]"
$ws.Range("B13").Value = "[]"
$ws.Range("C13").Value = "[]"
$ws.Range("D13").Value = "[]"
$ws.Range("E13").Value = "[]"
$ws.Range("F13").Value = "[]"
$ws.Range("B16").Value = "[]"
$ws.Range("C16").Value = "[]"
$ws.Range("D16").Value = "[]"
$ws.Range("E16").Value = "[]"
$ws.Range("F16").Value = "[]"
$ws.Range("B17").Value = "[]"
$ws.Range("C17").Value = "[]"
$ws.Range("D17").Value = "[]"
$ws.Range("E17").Value = "[61]"
$ws.Range("F17").Value = "[]"

# Sheet: best_seed
$ws = $wb.Worksheets.Item("best_seed")
$ws.Range("B2").Value = 29
$ws.Range("C2").Value = 29
$ws.Range("D2").Value = 7
$ws.Range("E2").Value = 31
$ws.Range("F2").Value = 13
$ws.Range("B3").Value = 19
$ws.Range("C3").Value = 19
$ws.Range("D3").Value = 17
$ws.Range("E3").Value = 37
$ws.Range("F3").Value = 59
$ws.Range("B4").Value = 23
$ws.Range("C4").Value = 11
$ws.Range("D4").Value = 17
$ws.Range("E4").Value = 37
$ws.Range("F4").Value = 2
$ws.Range("B5").Value = 71
$ws.Range("C5").Value = 13
$ws.Range("D5").Value = 7
$ws.Range("E5").Value = 37
$ws.Range("F5").Value = 3
$ws.Range("B6").Value = 29
$ws.Range("C6").Value = 47
$ws.Range("D6").Value = 43
$ws.Range("E6").Value = 37
$ws.Range("F6").Value = 13
$ws.Range("C7").Value = 31
$ws.Range("B9").Value = 19
$ws.Range("C9").Value = 17
$ws.Range("D9").Value = 7
$ws.Range("E9").Value = 37
$ws.Range("F9").Value = 29
$ws.Range("C10").Value = 37
$ws.Range("B11").Value = 41
$ws.Range("C11").Value = 71
$ws.Range("D11").Value = 53
$ws.Range("E11").Value = 37
$ws.Range("F11").Value = 31
$ws.Range("B12").Value = 19
$ws.Range("C12").Value = 13
$ws.Range("D12").Value = 3
$ws.Range("E12").Value = 11
$ws.Range("B13").Value = 29
$ws.Range("C13").Value = 59
$ws.Range("D13").Value = 61
$ws.Range("E13").Value = 23
$ws.Range("F13").Value = 59
$ws.Range("B16").Value = 71
$ws.Range("C16").Value = 13
$ws.Range("D16").Value = 19
$ws.Range("E16").Value = 37
$ws.Range("F16").Value = 3
$ws.Range("B17").Value = 19
$ws.Range("C17").Value = 19
$ws.Range("D17").Value = 61
$ws.Range("E17").Value = 37
$ws.Range("F17").Value = 13
